$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.1981799344354955
$ws.Range("J4").Value = 0.4570840386122355
$ws.Range("K4").Value = 0.267064564299616
$ws.Range("L4").Value = 2.636059852754527
